$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.004.90'
$ws.Range('E2').Value = '  -6.38%  '
$ws.Range('D3').Value = '3.388.39'
$ws.Range('E3').Value = '  -5.49%  '
$ws.Range('D4').Value = "'1.01"
$ws.Range('E4').Value = '  +1.05%  '
$ws.Range('D5').Value = "'385.63"
$ws.Range('E5').Value = '  -7.49%  '
$ws.Range('D6').Value = "'121.54"
$ws.Range('E6').Value = '  -5.97%  '
$ws.Range('D7').Value = '3.550.79'
$ws.Range('E7').Value = '  -0.68%  '
$ws.Range('D8').Value = "'0.576"
$ws.Range('E8').Value = '  -11.48%  '
$ws.Range('D9').Value = "'0.999"
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('D10').Value = "'0.657"
$ws.Range('E10').Value = '  -15.14%  '
$ws.Range('D11').Value = "'0.143"
$ws.Range('E11').Value = '  -19.41%  '
$ws.Range('D12').Value = "'0.0000293"
$ws.Range('E12').Value = '  -12.76%  '
$ws.Range('D13').Value = "'38.13"
$ws.Range('E13').Value = '  -10.22%  '
$ws.Range('D14').Value = '3.933.41'
$ws.Range('E14').Value = '  -5.47%  '
$ws.Range('D15').Value = "'8.99"
$ws.Range('E15').Value = '  -9.09%  '
$ws.Range('E16').Value = '  -2.97%  '
$ws.Range('D17').Value = '3.463.90'
$ws.Range('E17').Value = '  -3.72%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = "'18.27"
$ws.Range('E18').Value = '  -10.57%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = "'12.26"
$ws.Range('E19').Value = '  -0.25%  '
$ws.Range('D20').Value = '62.979.95'
$ws.Range('E20').Value = '  -6.27%  '
$ws.Range('D21').Value = "'1.01"
$ws.Range('E21').Value = '  -12.27%  '
$ws.Range('D22').Value = "'376.21"
$ws.Range('E22').Value = '  -16.66%  '
$ws.Range('D23').Value = "'13.47"
$ws.Range('E23').Value = '  +2.32%  '
$ws.Range('D24').Value = "'78.79"
$ws.Range('E24').Value = '  -11.83%  '
$ws.Range('D25').Value = "'2.76"
$ws.Range('E25').Value = '  -12.49%  '
$ws.Range('D26').Value = "'5.19"
$ws.Range('E26').Value = '  +6.30%  '
$ws.Range('D27').Value = "'32.69"
$ws.Range('E27').Value = '  -6.58%  '
$ws.Range('D28').Value = "'2.92"
$ws.Range('E28').Value = '  -13.13%  '
$ws.Range('D29').Value = "'8.61"
$ws.Range('E29').Value = '  -14.44%  '
$ws.Range('D30').Value = "'11.70"
$ws.Range('E30').Value = '  -5.58%  '
$ws.Range('D31').Value = "'2.59"
$ws.Range('E31').Value = '  -4.92%  '
$ws.Range('D32').Value = "'0.108"
$ws.Range('E32').Value = '  -8.46%  '
$ws.Range('D33').Value = "'6.46"
$ws.Range('E33').Value = '  -12.77%  '
$ws.Range('B34').Value = 'Dai'
$ws.Range('C34').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D34').Value = "'1.00"
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = "'0.145"
$ws.Range('E35').Value = '  -10.86%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = "'36.65"
$ws.Range('E36').Value = '  -11.62%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').Value = "'53.47"
$ws.Range('E37').Value = '  -5.72%  '
$ws.Range('D38').Value = "'0.0427"
$ws.Range('E38').Value = '  -13.67%  '
$ws.Range('D39').Value = "'0.990"
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('B40').Value = 'ThetaToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D40').Value = "'2.67"
$ws.Range('E40').Value = '  +15.16%  '
$ws.Range('B41').Value = 'EnergySwap'
$ws.Range('C41').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D41').Value = "'26.45"
$ws.Range('E41').Value = '  +23.40%  '
$ws.Range('B42').Value = 'Stellar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D42').Value = "'0.130"
$ws.Range('E42').Value = '  -11.27%  '
$ws.Range('B43').Value = 'Monero'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D43').Value = "'139.03"
$ws.Range('E43').Value = '  -6.68%  '
$ws.Range('D44').Value = "'2.98"
$ws.Range('E44').Value = '  +15.02%  '
$ws.Range('B45').Value = 'PEPE'
$ws.Range('C45').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D45').Value = '0.0₃0594'
$ws.Range('E45').Value = '  -20.31%  '
$ws.Range('D46').Value = "'1.91"
$ws.Range('E46').Value = '  -3.46%  '
$ws.Range('B47').Value = 'WEMIXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D47').Value = "'2.43"
$ws.Range('E47').Value = '  -11.88%  '
$ws.Range('D48').Value = "'2.99"
$ws.Range('E48').Value = '  -8.57%  '
$ws.Range('E49').Value = '  -9.23%  '
$ws.Range('D50').Value = "'2.60"
$ws.Range('E50').Value = '  -15.60%  '
$ws.Range('D51').Value = "'0.266"
$ws.Range('E51').Value = '  -15.98%  '
